# Applies the 230918 quiz-response update: appends 22 new survey rows
# (rows 491-512) to the single worksheet, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy the formatting (cell styles) of the last existing data row (490)
# columns A:L onto each of the new rows 491:512, since every new row reuses the
# exact same style indices as the rest of the data block (date/general/percent).
$ws.Range("A490:L490").Copy() | Out-Null
for ($r = 491; $r -le 512; $r++) {
    $ws.Range("A$r" + ":L$r").PasteSpecial(-4122) | Out-Null
}

# Step 2: the survey has two mutually exclusive trailing columns (M and N) -
# only one of them is populated per row. Copy the same text-cell format used
# elsewhere onto whichever single column (M or N) each new row actually uses,
# leaving the other column untouched (absent), exactly like the source rows.
$ws.Range("B490").Copy() | Out-Null
$ws.Range("N491").PasteSpecial(-4122) | Out-Null
$ws.Range("N492").PasteSpecial(-4122) | Out-Null
$ws.Range("M493").PasteSpecial(-4122) | Out-Null
$ws.Range("M494").PasteSpecial(-4122) | Out-Null
$ws.Range("N495").PasteSpecial(-4122) | Out-Null
$ws.Range("N496").PasteSpecial(-4122) | Out-Null
$ws.Range("M497").PasteSpecial(-4122) | Out-Null
$ws.Range("M498").PasteSpecial(-4122) | Out-Null
$ws.Range("M499").PasteSpecial(-4122) | Out-Null
$ws.Range("N500").PasteSpecial(-4122) | Out-Null
$ws.Range("M501").PasteSpecial(-4122) | Out-Null
$ws.Range("M502").PasteSpecial(-4122) | Out-Null
$ws.Range("N503").PasteSpecial(-4122) | Out-Null
$ws.Range("M504").PasteSpecial(-4122) | Out-Null
$ws.Range("N505").PasteSpecial(-4122) | Out-Null
$ws.Range("N506").PasteSpecial(-4122) | Out-Null
$ws.Range("M507").PasteSpecial(-4122) | Out-Null
$ws.Range("M508").PasteSpecial(-4122) | Out-Null
$ws.Range("M509").PasteSpecial(-4122) | Out-Null
$ws.Range("N510").PasteSpecial(-4122) | Out-Null
$ws.Range("N511").PasteSpecial(-4122) | Out-Null
$ws.Range("M512").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Step 3: populate the actual cell values for each new row.
# Row 491
$ws.Range("A491").Value = 45193.919750416666
$ws.Range("B491").Value = 'obj2107@gmail.com'
$ws.Range("C491").Value = '경제학과'
$ws.Range("D491").Value = 20173026
$ws.Range("E491").Value = '오병우'
$ws.Range("F491").Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Range("G491").Value = 0.3
$ws.Range("H491").Value = '6:4'
$ws.Range("I491").Value = '20분의 1'
$ws.Range("J491").Value = '20만호, 69만명'
$ws.Range("K491").Value = '경상'
$ws.Range("L491").Value = 'Black'
$ws.Range("N491").Value = '모름/무응답'

# Row 492
$ws.Range("A492").Value = 45193.92404366898
$ws.Range("B492").Value = 'happle0313@naver.com'
$ws.Range("C492").Value = '경영학과 '
$ws.Range("D492").Value = 20202945
$ws.Range("E492").Value = '김희경'
$ws.Range("F492").Value = '실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다.'
$ws.Range("G492").Value = 0.5
$ws.Range("H492").Value = '6:4'
$ws.Range("I492").Value = '20분의 1'
$ws.Range("J492").Value = '44만호, 153만명'
$ws.Range("K492").Value = '경상'
$ws.Range("L492").Value = 'Black'
$ws.Range("N492").Value = '모름/무응답'

# Row 493
$ws.Range("A493").Value = 45193.93739571759
$ws.Range("B493").Value = 'teslahan21@gmail.com'
$ws.Range("C493").Value = '소프트웨어학부'
$ws.Range("D493").Value = 20235276
$ws.Range("E493").Value = '한태웅'
$ws.Range("F493").Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Range("G493").Value = 0.1
$ws.Range("H493").Value = '6:4'
$ws.Range("I493").Value = '20분의 1'
$ws.Range("J493").Value = '20만호, 69만명'
$ws.Range("K493").Value = '충청'
$ws.Range("L493").Value = 'Red'
$ws.Range("M493").Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 494
$ws.Range("A494").Value = 45193.94430052083
$ws.Range("B494").Value = 'soojin020524@gmail.com'
$ws.Range("C494").Value = '체육학과'
$ws.Range("D494").Value = 20217130
$ws.Range("E494").Value = '노수진'
$ws.Range("F494").Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Range("G494").Value = 0.1
$ws.Range("H494").Value = '6:4'
$ws.Range("I494").Value = '20분의 1'
$ws.Range("J494").Value = '20만호, 69만명'
$ws.Range("K494").Value = '평안'
$ws.Range("L494").Value = 'Red'
$ws.Range("M494").Value = '모름/무응답'

# Row 495
$ws.Range("A495").Value = 45193.945696863426
$ws.Range("B495").Value = 'ryan.jin1005@gmail.com'
$ws.Range("C495").Value = '금융재무학과'
$ws.Range("D495").Value = 20192999
$ws.Range("E495").Value = '진승재'
$ws.Range("F495").Value = '실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다.'
$ws.Range("G495").Value = 0.3
$ws.Range("H495").Value = '5:5'
$ws.Range("I495").Value = '10분의 1'
$ws.Range("J495").Value = '15만호,  32만명'
$ws.Range("K495").Value = '평안'
$ws.Range("L495").Value = 'Black'
$ws.Range("N495").Value = '모름/무응답'

# Row 496
$ws.Range("A496").Value = 45193.95218641203
$ws.Range("B496").Value = 'dlgusdnr4580@naver.com'
$ws.Range("C496").Value = '러시아학거ㅏ'
$ws.Range("D496").Value = 20101722
$ws.Range("E496").Value = '이현욱'
$ws.Range("F496").Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Range("G496").Value = 0.9
$ws.Range("H496").Value = '5:5'
$ws.Range("I496").Value = '10분의 1'
$ws.Range("J496").Value = '15만호,  32만명'
$ws.Range("K496").Value = '충청'
$ws.Range("L496").Value = 'Black'
$ws.Range("N496").Value = '모름/무응답'

# Row 497
$ws.Range("A497").Value = 45193.955362037035
$ws.Range("B497").Value = 'bin85288@daum.net'
$ws.Range("C497").Value = '인문학부'
$ws.Range("D497").Value = 20231016
$ws.Range("E497").Value = '김수빈'
$ws.Range("F497").Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Range("G497").Value = 0.9
$ws.Range("H497").Value = '6:4'
$ws.Range("I497").Value = '20분의 1'
$ws.Range("J497").Value = '44만호, 153만명'
$ws.Range("K497").Value = '충청'
$ws.Range("L497").Value = 'Red'
$ws.Range("M497").Value = '반대한다.'

# Row 498
$ws.Range("A498").Value = 45193.956515046295
$ws.Range("B498").Value = 'kimlee0411@naver.com'
$ws.Range("C498").Value = '사회복학부'
$ws.Range("D498").Value = 20232314
$ws.Range("E498").Value = '김수빈'
$ws.Range("F498").Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Range("G498").Value = 0.9
$ws.Range("H498").Value = '4:6'
$ws.Range("I498").Value = '15분의 1'
$ws.Range("J498").Value = '15만호,  32만명'
$ws.Range("K498").Value = '경상'
$ws.Range("L498").Value = 'Red'
$ws.Range("M498").Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 499
$ws.Range("A499").Value = 45193.96031199074
$ws.Range("B499").Value = 'ktm4145@gmail.com'
$ws.Range("C499").Value = '바이오메디컬학과'
$ws.Range("D499").Value = 20223608
$ws.Range("E499").Value = '김태민'
$ws.Range("F499").Value = '실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다.'
$ws.Range("G499").Value = 0.7
$ws.Range("H499").Value = '5:5'
$ws.Range("I499").Value = '20분의 1'
$ws.Range("J499").Value = '20만호, 69만명'
$ws.Range("K499").Value = '전라'
$ws.Range("L499").Value = 'Red'
$ws.Range("M499").Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 500
$ws.Range("A500").Value = 45193.962819548615
$ws.Range("B500").Value = 'enjoy9675@gmail.com'
$ws.Range("C500").Value = '바이오메디컬'
$ws.Range("D500").Value = 20213801
$ws.Range("E500").Value = '곽아영'
$ws.Range("F500").Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Range("G500").Value = 0.1
$ws.Range("H500").Value = '6:4'
$ws.Range("I500").Value = '20분의 1'
$ws.Range("J500").Value = '20만호, 69만명'
$ws.Range("K500").Value = '충청'
$ws.Range("L500").Value = 'Black'
$ws.Range("N500").Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 501
$ws.Range("A501").Value = 45193.96440173611
$ws.Range("B501").Value = 'heeyeon_02@naver.com'
$ws.Range("C501").Value = '중국학과'
$ws.Range("D501").Value = 20231536
$ws.Range("E501").Value = '이희연'
$ws.Range("F501").Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Range("G501").Value = 0.3
$ws.Range("H501").Value = '4:6'
$ws.Range("I501").Value = '15분의 1'
$ws.Range("J501").Value = '44만호, 153만명'
$ws.Range("K501").Value = '경기'
$ws.Range("L501").Value = 'Red'
$ws.Range("M501").Value = '모름/무응답'

# Row 502
$ws.Range("A502").Value = 45193.96647329861
$ws.Range("B502").Value = 'b1a46317@naver.com'
$ws.Range("C502").Value = '경영학과'
$ws.Range("D502").Value = 20203020
$ws.Range("E502").Value = '이정수'
$ws.Range("F502").Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Range("G502").Value = 0.7
$ws.Range("H502").Value = '5:5'
$ws.Range("I502").Value = '15분의 1'
$ws.Range("J502").Value = '20만호, 69만명'
$ws.Range("K502").Value = '평안'
$ws.Range("L502").Value = 'Red'
$ws.Range("M502").Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 503
$ws.Range("A503").Value = 45193.96754217593
$ws.Range("B503").Value = 'sherisim@naver.com'
$ws.Range("C503").Value = '인공지능융합학부'
$ws.Range("D503").Value = 20236741
$ws.Range("E503").Value = '심시원'
$ws.Range("F503").Value = '실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다.'
$ws.Range("G503").Value = 0.3
$ws.Range("H503").Value = '6:4'
$ws.Range("I503").Value = '20분의 1'
$ws.Range("J503").Value = '20만호, 69만명'
$ws.Range("K503").Value = '전라'
$ws.Range("L503").Value = 'Black'
$ws.Range("N503").Value = '찬성한다.'

# Row 504
$ws.Range("A504").Value = 45193.97474690972
$ws.Range("B504").Value = 'yoogyeonggg@naver.com'
$ws.Range("C504").Value = '영어영문학과'
$ws.Range("D504").Value = 20231238
$ws.Range("E504").Value = '최유경'
$ws.Range("F504").Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Range("G504").Value = 0.7
$ws.Range("H504").Value = '6:4'
$ws.Range("I504").Value = '20분의 1'
$ws.Range("J504").Value = '20만호, 69만명'
$ws.Range("K504").Value = '경상'
$ws.Range("L504").Value = 'Red'
$ws.Range("M504").Value = '모름/무응답'

# Row 505
$ws.Range("A505").Value = 45193.975047581014
$ws.Range("B505").Value = 'scw0922@naver.com'
$ws.Range("C505").Value = '간호학과'
$ws.Range("D505").Value = 20236256
$ws.Range("E505").Value = '신채원'
$ws.Range("F505").Value = '과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다.'
$ws.Range("G505").Value = 0.1
$ws.Range("H505").Value = '7:3'
$ws.Range("I505").Value = '10분의 1'
$ws.Range("J505").Value = '15만호,  32만명'
$ws.Range("K505").Value = '경기'
$ws.Range("L505").Value = 'Black'
$ws.Range("N505").Value = '모름/무응답'

# Row 506
$ws.Range("A506").Value = 45193.97682863426
$ws.Range("B506").Value = 'thdrkdud456456@gmail.com'
$ws.Range("C506").Value = '경영학과'
$ws.Range("D506").Value = 20222968
$ws.Range("E506").Value = '송가영'
$ws.Range("F506").Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Range("G506").Value = 0.5
$ws.Range("H506").Value = '3:7'
$ws.Range("I506").Value = '10분의 1'
$ws.Range("J506").Value = '130만호, 5백만명'
$ws.Range("K506").Value = '전라'
$ws.Range("L506").Value = 'Black'
$ws.Range("N506").Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 507
$ws.Range("A507").Value = 45193.97683835648
$ws.Range("B507").Value = 'kxjenlee@naver.com'
$ws.Range("C507").Value = '글로벌비즈니스'
$ws.Range("D507").Value = 20226417
$ws.Range("E507").Value = '이제인'
$ws.Range("F507").Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Range("G507").Value = 0.7
$ws.Range("H507").Value = '6:4'
$ws.Range("I507").Value = '20분의 1'
$ws.Range("J507").Value = '44만호, 153만명'
$ws.Range("K507").Value = '전라'
$ws.Range("L507").Value = 'Red'
$ws.Range("M507").Value = '모름/무응답'

# Row 508
$ws.Range("A508").Value = 45193.982407013886
$ws.Range("B508").Value = 'lju5422@naver.com'
$ws.Range("C508").Value = '사회복지학부'
$ws.Range("D508").Value = 20192348
$ws.Range("E508").Value = '임준섭'
$ws.Range("F508").Value = '‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다.'
$ws.Range("G508").Value = 0.9
$ws.Range("H508").Value = '5:5'
$ws.Range("I508").Value = '30분의 1'
$ws.Range("J508").Value = '44만호, 153만명'
$ws.Range("K508").Value = '경기'
$ws.Range("L508").Value = 'Red'
$ws.Range("M508").Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 509
$ws.Range("A509").Value = 45193.98299708334
$ws.Range("B509").Value = 'choe0119@gmail.com'
$ws.Range("C509").Value = '의예과'
$ws.Range("D509").Value = 20226176
$ws.Range("E509").Value = '최태웅'
$ws.Range("F509").Value = '과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다.'
$ws.Range("G509").Value = 0.3
$ws.Range("H509").Value = '6:4'
$ws.Range("I509").Value = '10분의 1'
$ws.Range("J509").Value = '20만호, 69만명'
$ws.Range("K509").Value = '전라'
$ws.Range("L509").Value = 'Red'
$ws.Range("M509").Value = '반대한다.'

# Row 510
$ws.Range("A510").Value = 45193.99949899306
$ws.Range("B510").Value = 'hanhj0223@gmail.com'
$ws.Range("C510").Value = '언어청각학부'
$ws.Range("D510").Value = 20233966
$ws.Range("E510").Value = '한형준'
$ws.Range("F510").Value = '과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다.'
$ws.Range("G510").Value = 0.5
$ws.Range("H510").Value = '4:6'
$ws.Range("I510").Value = '15분의 1'
$ws.Range("J510").Value = '15만호,  32만명'
$ws.Range("K510").Value = '평안'
$ws.Range("L510").Value = 'Black'
$ws.Range("N510").Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 511
$ws.Range("A511").Value = 45194.00363949074
$ws.Range("B511").Value = 'chi9605@naver.com'
$ws.Range("C511").Value = '광고홍보학과'
$ws.Range("D511").Value = 20202633
$ws.Range("E511").Value = '이채연'
$ws.Range("F511").Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Range("G511").Value = 0.1
$ws.Range("H511").Value = '6:4'
$ws.Range("I511").Value = '20분의 1'
$ws.Range("J511").Value = '20만호, 69만명'
$ws.Range("K511").Value = '충청'
$ws.Range("L511").Value = 'Black'
$ws.Range("N511").Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 512
$ws.Range("A512").Value = 45194.004796284724
$ws.Range("B512").Value = 'lian_ycm@naver.com'
$ws.Range("C512").Value = '광고홍보학과'
$ws.Range("D512").Value = 20202622
$ws.Range("E512").Value = '유채민'
$ws.Range("F512").Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Range("G512").Value = 0.1
$ws.Range("H512").Value = '6:4'
$ws.Range("I512").Value = '20분의 1'
$ws.Range("J512").Value = '20만호, 69만명'
$ws.Range("K512").Value = '충청'
$ws.Range("L512").Value = 'Red'
$ws.Range("M512").Value = '반대한다.'
